$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the two phone numbers to international (E.164-ish) format:
# "050-7676706" -> "972507676706", "058-6208430" -> "972586208430"
$ws.Range("B2").Value = "972507676706"
$ws.Range("B3").Value = "972586208430"

# The longer digit strings no longer fit the old best-fit width, so widen
# column B (stored width needs to land on exactly 13).
$ws.Columns.Item(2).ColumnWidth = 12.3

# Move/save the active selection to E8
$ws.Range("E8").Select()
